# "added slide for references"
# Adds a new slide (6th slide / sldId 262) containing two text boxes:
#   1) a "References" title textbox
#   2) a textbox with a hyperlinked "Feature maps" line and a bulleted
#      hyperlinked URL line, followed by a blank trailing paragraph.

$p = $ppt.ActivePresentation

# EMU -> point conversion helper (1 pt = 12700 EMU)
function EMU($v) { return $v / 12700.0 }

# Insert the new slide at the end, using the same blank layout (index 7,
# "Blank") that the rest of the deck uses. ppLayoutBlank = 12.
$newIndex = $p.Slides.Count + 1
$slide = $p.Slides.Add($newIndex, 12)

# ---- Shape 2: "References" title textbox ------------------------------
$title = $slide.Shapes.AddTextbox(1, (EMU 3941805), (EMU 543697), (EMU 3323968), (EMU 369332))
$title.Fill.Visible = $false
$title.TextFrame.WordWrap = $true
$title.TextFrame.AutoSize = 1

$titleRange = $title.TextFrame.TextRange
$titleRange.Text = "References"

# ---- Shape 3: "Feature maps" + link textbox ----------------------------
$body = $slide.Shapes.AddTextbox(1, (EMU 729049), (EMU 1186249), (EMU 5474043), (EMU 923330))
$body.Fill.Visible = $false
$body.TextFrame.WordWrap = $true
$body.TextFrame.AutoSize = 1

$linkUrl = "https://arxiv.org/pdf/1507.02313.pdf"

# Build up the paragraphs one at a time, applying the hyperlink while the
# textbox still holds a single paragraph and then growing it with
# InsertAfter; the engine's hyperlink assignment otherwise clobbers the
# run-level language attribute once a shape already has >1 paragraph.
$bodyRange = $body.TextFrame.TextRange
$bodyRange.Text = "Feature maps"
$bodyRange.ActionSettings(3).Hyperlink.Address = $linkUrl

$bodyRange.InsertAfter("`r" + $linkUrl)

# Trailing blank paragraph: add a throwaway character (inherits the
# hyperlink formatting from InsertAfter), strip its hyperlink, then
# delete it again - leaves a clean empty paragraph behind.
$bodyRange.InsertAfter("`rX")
$fullRange = $body.TextFrame.TextRange
$lastChar = $fullRange.Characters($fullRange.Length, 1)
$lastChar.ActionSettings(3).Hyperlink.Address = ""
$lastChar.Text = ""

# Paragraph 2 (the URL line) gets a bullet.
$para2 = $fullRange.Paragraphs(2, 1)
$para2.ParagraphFormat.Bullet.Visible = $true
$para2.ParagraphFormat.Bullet.Character = 8226
$para2.ParagraphFormat.Bullet.Font.Name = "Arial"

Write-Host "Added references slide; slide count now:" $p.Slides.Count
